$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The previous "last" payment row (row 19, phone 79174445) gets its phone
# number normalized from text to a real number now that it's no longer the
# newest entry.
$ws.Range("A19").Value = 79174445

# Append the new payment as row 20. The phone number is written as text
# (leading apostrophe forces text, matching how new rows are first recorded)
# while the numeric fields are written as real numbers. The empty
# discount/method placeholder columns (B, F) are written as empty text,
# matching the blank-but-text cells used elsewhere in the sheet.
$ws.Range("A20").Value = "'79174445"
$ws.Cells.Item(20, 2).Value = "'"
$ws.Range("C20").Value = "Cash"
$ws.Range("D20").Value = "2025-08-18T09:08:11"
$ws.Range("E20").Value = 20
$ws.Cells.Item(20, 6).Value = "'"
$ws.Range("G20").Value = 0
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 20
